$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation row is inserted right before the current row 750
# (2026/12/29, 火, 13, 201), pushing it and every row below it down by
# one. The newly inserted row 750 holds the "missing" 2026/02/01 reading
# that belongs right after the existing 2026/02/01 03:00 row (row 749).
$ws.Rows(750).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/02/01"), not
# real Excel date serials. Force the cell to text *before* writing so the
# "2026/02/01" string isn't auto-converted into a date value, then reset
# the style back to Normal so no extra (quote-prefix) formatting sticks
# around on the cell - matching the unstyled cells around it.
$ws.Cells.Item(750, 1).NumberFormat = "@"
$ws.Cells.Item(750, 1).Value = "2026/02/01"
$ws.Cells.Item(750, 1).Style = "Normal"

$ws.Cells.Item(750, 2).Value = "日"
$ws.Cells.Item(750, 3).Value = 7
$ws.Cells.Item(750, 4).Value = 201
